$wb = $excel.ActiveWorkbook

# New execution timestamps for the FEINmismatch sheet (B2:B30)
$feinMismatchTimestamps = @(
    "Fri Oct 25 11:21:32 EDT 2024",
    "Fri Oct 25 11:21:46 EDT 2024",
    "Fri Oct 25 11:21:58 EDT 2024",
    "Fri Oct 25 11:22:09 EDT 2024",
    "Fri Oct 25 11:22:21 EDT 2024",
    "Fri Oct 25 11:22:33 EDT 2024",
    "Fri Oct 25 11:22:44 EDT 2024",
    "Fri Oct 25 11:22:57 EDT 2024",
    "Fri Oct 25 11:23:08 EDT 2024",
    "Fri Oct 25 11:23:19 EDT 2024",
    "Fri Oct 25 11:23:31 EDT 2024",
    "Fri Oct 25 11:23:42 EDT 2024",
    "Fri Oct 25 11:23:54 EDT 2024",
    "Fri Oct 25 11:24:05 EDT 2024",
    "Fri Oct 25 11:24:16 EDT 2024",
    "Fri Oct 25 11:24:28 EDT 2024",
    "Fri Oct 25 11:24:40 EDT 2024",
    "Fri Oct 25 11:24:51 EDT 2024",
    "Fri Oct 25 11:25:02 EDT 2024",
    "Fri Oct 25 11:25:14 EDT 2024",
    "Fri Oct 25 11:25:25 EDT 2024",
    "Fri Oct 25 11:25:37 EDT 2024",
    "Fri Oct 25 11:25:49 EDT 2024",
    "Fri Oct 25 11:26:01 EDT 2024",
    "Fri Oct 25 11:26:12 EDT 2024",
    "Fri Oct 25 11:26:23 EDT 2024",
    "Fri Oct 25 11:26:34 EDT 2024",
    "Fri Oct 25 11:26:46 EDT 2024",
    "Fri Oct 25 11:26:57 EDT 2024"
)

# New execution timestamps for the FEINSSNmismatch sheet (B2:B19)
$feinSsnMismatchTimestamps = @(
    "Fri Oct 25 11:27:09 EDT 2024",
    "Fri Oct 25 11:27:20 EDT 2024",
    "Fri Oct 25 11:27:31 EDT 2024",
    "Fri Oct 25 11:27:42 EDT 2024",
    "Fri Oct 25 11:27:53 EDT 2024",
    "Fri Oct 25 11:28:04 EDT 2024",
    "Fri Oct 25 11:28:15 EDT 2024",
    "Fri Oct 25 11:28:26 EDT 2024",
    "Fri Oct 25 11:28:37 EDT 2024",
    "Fri Oct 25 11:28:48 EDT 2024",
    "Fri Oct 25 11:28:59 EDT 2024",
    "Fri Oct 25 11:29:11 EDT 2024",
    "Fri Oct 25 11:29:22 EDT 2024",
    "Fri Oct 25 11:29:33 EDT 2024",
    "Fri Oct 25 11:29:44 EDT 2024",
    "Fri Oct 25 11:29:55 EDT 2024",
    "Fri Oct 25 11:30:07 EDT 2024",
    "Fri Oct 25 11:30:18 EDT 2024"
)

$wsFeinMismatch = $wb.Worksheets.Item("FEINmismatch")
for ($i = 0; $i -lt $feinMismatchTimestamps.Length; $i++) {
    $wsFeinMismatch.Cells.Item($i + 2, 2).Value = $feinMismatchTimestamps[$i]
}

$wsFeinSsnMismatch = $wb.Worksheets.Item("FEINSSNmismatch")
for ($i = 0; $i -lt $feinSsnMismatchTimestamps.Length; $i++) {
    $wsFeinSsnMismatch.Cells.Item($i + 2, 2).Value = $feinSsnMismatchTimestamps[$i]
}
